$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 96; this pushes the existing rows 96-148
# down to 99-151 and grows the used range to A1:R151 automatically.
$ws.Rows.Item(96).Resize(3).Insert()

# --- New row 96 (Banquete) ---
$ws.Range("A96").Value = 6
$ws.Range("B96").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C96").Value = "Metropolitana"
$ws.Range("D96").Value = 44875
$ws.Range("E96").Value = 13
$ws.Range("F96").Value = 300000000
$ws.Range("G96").Value = "Espárragos"
$ws.Range("H96").Value = "Sin especificar"
$ws.Range("I96").Value = "Banquete"
$ws.Range("J96").Value = 1400
$ws.Range("K96").Value = 1400
$ws.Range("L96").Value = 1500
$ws.Range("M96").Value = 1454
$ws.Range("N96").Value = "`$/kilo"
$ws.Range("O96").Value = "Provincia de Linares"
$ws.Range("P96").Value = 1454
$ws.Range("Q96").Value = 1
$ws.Range("R96").Value = "Hortaliza"

# --- New row 97 (Primera) ---
$ws.Range("A97").Value = 6
$ws.Range("B97").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C97").Value = "Metropolitana"
$ws.Range("D97").Value = 44875
$ws.Range("E97").Value = 13
$ws.Range("F97").Value = 300000000
$ws.Range("G97").Value = "Espárragos"
$ws.Range("H97").Value = "Sin especificar"
$ws.Range("I97").Value = "Primera"
$ws.Range("J97").Value = 1550
$ws.Range("K97").Value = 1200
$ws.Range("L97").Value = 1300
$ws.Range("M97").Value = 1252
$ws.Range("N97").Value = "`$/kilo"
$ws.Range("O97").Value = "Provincia de Linares"
$ws.Range("P97").Value = 1252
$ws.Range("Q97").Value = 1
$ws.Range("R97").Value = "Hortaliza"

# --- New row 98 (Segunda) ---
$ws.Range("A98").Value = 6
$ws.Range("B98").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C98").Value = "Metropolitana"
$ws.Range("D98").Value = 44875
$ws.Range("E98").Value = 13
$ws.Range("F98").Value = 300000000
$ws.Range("G98").Value = "Espárragos"
$ws.Range("H98").Value = "Sin especificar"
$ws.Range("I98").Value = "Segunda"
$ws.Range("J98").Value = 1150
$ws.Range("K98").Value = 1000
$ws.Range("L98").Value = 1100
$ws.Range("M98").Value = 1052
$ws.Range("N98").Value = "`$/kilo"
$ws.Range("O98").Value = "Provincia de Linares"
$ws.Range("P98").Value = 1052
$ws.Range("Q98").Value = 1
$ws.Range("R98").Value = "Hortaliza"
